$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5276
$ws.Range("I15").Value = 5276
$ws.Range("K15").Value = 15828
$ws.Range("M15").Value = -15659
$ws.Range("H109").Value = 35468.4
$ws.Range("J109").Value = 35468.4
$ws.Range("L109").Value = 35468.4
$ws.Range("N109").Value = -38242.4
$ws.Range("H132").Value = 20661.469
$ws.Range("I132").Value = 2947.3257
$ws.Range("J132").Value = 147612.83
$ws.Range("K132").Value = 8841.9771
$ws.Range("L132").Value = 442838.49
$ws.Range("M132").Value = -6311.9771
$ws.Range("N132").Value = -447898.49
$ws.Range("H133").Value = 39337.785
$ws.Range("J133").Value = 39337.785
$ws.Range("L133").Value = 39337.785
$ws.Range("N133").Value = -49457.785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12495.389
$ws.Range("I32").Value = 11516.892
$ws.Range("K32").Value = 11516.892
$ws.Range("M32").Value = -11229.892
$ws.Range("H45").Value = 2313.48
$ws.Range("I45").Value = 1930.2858
$ws.Range("J45").Value = 2801.182
$ws.Range("K45").Value = 1930.2858
$ws.Range("L45").Value = 2801.182
$ws.Range("M45").Value = -1553.2858
$ws.Range("N45").Value = -3555.182
$ws.Range("H61").Value = 1581.2433
$ws.Range("I61").Value = 1323.8529
$ws.Range("K61").Value = 1323.8529
$ws.Range("M61").Value = -1111.8529
$ws.Range("H74").Value = 1487.6545
$ws.Range("I74").Value = 1487.0513
$ws.Range("J74").Value = 1489.125
$ws.Range("K74").Value = 1487.0513
$ws.Range("L74").Value = 1489.125
$ws.Range("M74").Value = -613.0513000000001
$ws.Range("N74").Value = -3237.125
$ws.Range("H77").Value = 1487.6545
$ws.Range("I77").Value = 1487.0513
$ws.Range("J77").Value = 1489.125
$ws.Range("K77").Value = 7435.2565
$ws.Range("L77").Value = 7445.625
$ws.Range("M77").Value = -3067.2565
$ws.Range("N77").Value = -16181.625
$ws.Range("H102").Value = 12660.952
$ws.Range("I102").Value = 2606.1538
$ws.Range("J102").Value = 29000
$ws.Range("K102").Value = 2606.1538
$ws.Range("L102").Value = 29000
$ws.Range("M102").Value = -984.1538
$ws.Range("N102").Value = -32244
$ws.Range("H110").Value = 1497.8387
$ws.Range("I110").Value = 1501.32
$ws.Range("K110").Value = 1501.32
$ws.Range("M110").Value = 543.6800000000001
$ws.Range("H117").Value = 42573.4
$ws.Range("J117").Value = 42573.4
$ws.Range("L117").Value = 42573.4
$ws.Range("N117").Value = -51751.4
$ws.Range("H122").Value = 1622.9574
$ws.Range("I122").Value = 1616.2195
$ws.Range("J122").Value = 1669
$ws.Range("K122").Value = 4848.6585
$ws.Range("L122").Value = 5007
$ws.Range("M122").Value = -2398.6585
$ws.Range("N122").Value = -9907
$ws.Range("H132").Value = 12502208
$ws.Range("I132").Value = 21740546
$ws.Range("J132").Value = 3278.7646
$ws.Range("K132").Value = 65221638
$ws.Range("L132").Value = 9836.293799999999
$ws.Range("M132").Value = -65219108
$ws.Range("N132").Value = -14896.2938
$ws.Range("H136").Value = 1581.2433
$ws.Range("I136").Value = 1323.8529
$ws.Range("K136").Value = 3971.5587
$ws.Range("M136").Value = -1421.5587
$ws.Range("H139").Value = 42762.1
$ws.Range("J139").Value = 42762.1
$ws.Range("L139").Value = 42762.1
$ws.Range("N139").Value = -53042.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1850.2444
$ws.Range("I105").Value = 1621.6111
$ws.Range("J105").Value = 2002.6666
$ws.Range("K105").Value = 1621.6111
$ws.Range("L105").Value = 2002.6666
$ws.Range("M105").Value = 125.3888999999999
$ws.Range("N105").Value = -5496.6666
$ws.Range("H107").Value = 1970.1082
$ws.Range("I107").Value = 1897.2084
$ws.Range("J107").Value = 2104.6924
$ws.Range("K107").Value = 1897.2084
$ws.Range("L107").Value = 2104.6924
$ws.Range("M107").Value = 22.79160000000002
$ws.Range("N107").Value = -5944.6924
$ws.Range("H112").Value = 46733
$ws.Range("J112").Value = 46733
$ws.Range("L112").Value = 46733
$ws.Range("N112").Value = -49687
$ws.Range("H116").Value = 44267.332
$ws.Range("J116").Value = 44267.332
$ws.Range("L116").Value = 44267.332
$ws.Range("N116").Value = -53445.332
$ws.Range("H132").Value = 44476.668
$ws.Range("J132").Value = 44476.668
$ws.Range("L132").Value = 44476.668
$ws.Range("N132").Value = -54596.668
$ws.Range("H134").Value = 3148.611
$ws.Range("I134").Value = 1621.5714
$ws.Range("J134").Value = 3777.392
$ws.Range("K134").Value = 4864.7142
$ws.Range("L134").Value = 11332.176
$ws.Range("M134").Value = -2329.7142
$ws.Range("N134").Value = -16402.176

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3016.77
$ws.Range("I31").Value = 1281.0869
$ws.Range("J31").Value = 3535.2207
$ws.Range("K31").Value = 1281.0869
$ws.Range("L31").Value = 3535.2207
$ws.Range("M31").Value = -986.0869
$ws.Range("N31").Value = -4125.2207
$ws.Range("H34").Value = 3016.77
$ws.Range("I34").Value = 1281.0869
$ws.Range("J34").Value = 3535.2207
$ws.Range("K34").Value = 1281.0869
$ws.Range("L34").Value = 3535.2207
$ws.Range("M34").Value = -1079.0869
$ws.Range("N34").Value = -3939.2207
$ws.Range("H58").Value = 1092.3539
$ws.Range("I58").Value = 837.7778
$ws.Range("J58").Value = 2342.0908
$ws.Range("K58").Value = 837.7778
$ws.Range("L58").Value = 2342.0908
$ws.Range("M58").Value = -634.7778
$ws.Range("N58").Value = -2748.0908
$ws.Range("H116").Value = 48489
$ws.Range("J116").Value = 48489
$ws.Range("L116").Value = 48489
$ws.Range("N116").Value = -57667
$ws.Range("H122").Value = 134733.33
$ws.Range("I122").Value = 301075
$ws.Range("J122").Value = 1660
$ws.Range("K122").Value = 903225
$ws.Range("L122").Value = 4980
$ws.Range("M122").Value = -900775
$ws.Range("N122").Value = -9880
$ws.Range("H134").Value = 2529.5264
$ws.Range("I134").Value = 1874.7273
$ws.Range("K134").Value = 5624.1819
$ws.Range("M134").Value = -3089.1819
$ws.Range("H136").Value = 1092.3539
$ws.Range("I136").Value = 837.7778
$ws.Range("J136").Value = 2342.0908
$ws.Range("K136").Value = 2513.3334
$ws.Range("L136").Value = 7026.2724
$ws.Range("M136").Value = 36.66660000000002
$ws.Range("N136").Value = -12126.2724
$ws.Range("H137").Value = 62299.5
$ws.Range("J137").Value = 62299.5
$ws.Range("L137").Value = 62299.5
$ws.Range("N137").Value = -72499.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7343.7144
$ws.Range("I97").Value = 3479
$ws.Range("J97").Value = 17005.5
$ws.Range("K97").Value = 3479
$ws.Range("L97").Value = 17005.5
$ws.Range("M97").Value = -2983
$ws.Range("N97").Value = -17997.5
$ws.Range("H102").Value = 1430.3182
$ws.Range("I102").Value = 1602.8667
$ws.Range("J102").Value = 1060.5714
$ws.Range("K102").Value = 1602.8667
$ws.Range("L102").Value = 1060.5714
$ws.Range("M102").Value = 19.13329999999996
$ws.Range("N102").Value = -4304.5714
$ws.Range("H113").Value = 1341.2273
$ws.Range("J113").Value = 1225.4286
$ws.Range("L113").Value = 1225.4286
$ws.Range("N113").Value = -5565.4286
$ws.Range("H114").Value = 38532.8
$ws.Range("J114").Value = 38532.8
$ws.Range("L114").Value = 38532.8
$ws.Range("N114").Value = -47210.8
$ws.Range("H122").Value = 1321.5385
$ws.Range("I122").Value = 1420
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 4260
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -1810
$ws.Range("N122").Value = -8200
$ws.Range("H135").Value = 54483
$ws.Range("J135").Value = 54483
$ws.Range("L135").Value = 54483
$ws.Range("N135").Value = -64623
$ws.Range("H138").Value = 36500
$ws.Range("J138").Value = 36500
$ws.Range("L138").Value = 36500
$ws.Range("N138").Value = -46780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2885
$ws.Range("I16").Value = 2664.7646
$ws.Range("K16").Value = 2664.7646
$ws.Range("M16").Value = -2494.7646
$ws.Range("H81").Value = 32181
$ws.Range("J81").Value = 32181
$ws.Range("L81").Value = 32181
$ws.Range("N81").Value = -34177
$ws.Range("H84").Value = 32181
$ws.Range("J84").Value = 32181
$ws.Range("L84").Value = 96543
$ws.Range("N84").Value = -106527
$ws.Range("H111").Value = 44379
$ws.Range("J111").Value = 44379
$ws.Range("L111").Value = 44379
$ws.Range("N111").Value = -52559

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 62044
$ws.Range("J46").Value = 62044
$ws.Range("L46").Value = 62044
$ws.Range("N46").Value = -62506
$ws.Range("H117").Value = 39170.168
$ws.Range("J117").Value = 39170.168
$ws.Range("L117").Value = 39170.168
$ws.Range("N117").Value = -48348.168
$ws.Range("H119").Value = 48674
$ws.Range("J119").Value = 48674
$ws.Range("L119").Value = 48674
$ws.Range("N119").Value = -58350
$ws.Range("H121").Value = 43412
$ws.Range("J121").Value = 43412
$ws.Range("L121").Value = 43412
$ws.Range("N121").Value = -46906
$ws.Range("H132").Value = 1839.317
$ws.Range("I132").Value = 1379.3793
$ws.Range("J132").Value = 2950.8333
$ws.Range("K132").Value = 4138.1379
$ws.Range("L132").Value = 8852.499899999999
$ws.Range("M132").Value = -1608.1379
$ws.Range("N132").Value = -13912.4999
$ws.Range("H134").Value = 62044
$ws.Range("J134").Value = 62044
$ws.Range("L134").Value = 186132
$ws.Range("N134").Value = -191202
$ws.Range("H138").Value = 45525
$ws.Range("J138").Value = 45525
$ws.Range("L138").Value = 45525
$ws.Range("N138").Value = -55805
